$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '37.481.76'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +2.67%  '
$ws.Range('E2').Style = "Normal"

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.067.75'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +3.55%  '
$ws.Range('E3').Style = "Normal"

# Row 4
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('E4').Style = "Normal"

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.95'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('E5').Style = "Normal"

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.618'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +2.87%  '
$ws.Range('E6').Style = "Normal"

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.29'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +6.38%  '
$ws.Range('E7').Style = "Normal"

# Row 8
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E8').Style = "Normal"

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.383'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +3.48%  '
$ws.Range('E9').Style = "Normal"

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '58.89'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E10').Style = "Normal"

# Row 11
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.93%  '
$ws.Range('E11').Style = "Normal"

# Row 12
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +3.41%  '
$ws.Range('E12').Style = "Normal"

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.372.74'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +3.52%  '
$ws.Range('E13').Style = "Normal"

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '14.52'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.39%  '
$ws.Range('E14').Style = "Normal"

# Row 15
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +4.36%  '
$ws.Range('E15').Style = "Normal"

# Row 16
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +2.84%  '
$ws.Range('E16').Style = "Normal"

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.19'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.26%  '
$ws.Range('E17').Style = "Normal"

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.064.64'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +2.80%  '
$ws.Range('E18').Style = "Normal"

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '37.662.84'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('E19').Style = "Normal"

# Row 20
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +17.58%  '
$ws.Range('E20').Style = "Normal"

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '70.14'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +3.61%  '
$ws.Range('E21').Style = "Normal"

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0816'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.64%  '
$ws.Range('E22').Style = "Normal"

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '226.78'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('E23').Style = "Normal"

# Row 24
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E24').Style = "Normal"

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.44'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.97%  '
$ws.Range('E25').Style = "Normal"

# Row 26
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.42%  '
$ws.Range('E26').Style = "Normal"

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '166.53'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.29%  '
$ws.Range('E27').Style = "Normal"

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.50'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +8.60%  '
$ws.Range('E28').Style = "Normal"

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.90'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.58%  '
$ws.Range('E29').Style = "Normal"

# Row 30
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.24'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +2.49%  '
$ws.Range('E30').Style = "Normal"

# Row 31
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.128'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.51%  '
$ws.Range('E31').Style = "Normal"

# Row 32
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.79%  '
$ws.Range('E32').Style = "Normal"

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.53'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +3.59%  '
$ws.Range('E33').Style = "Normal"

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0622'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.29%  '
$ws.Range('E34').Style = "Normal"

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.56'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +8.80%  '
$ws.Range('E35').Style = "Normal"

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.58'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +7.76%  '
$ws.Range('E36').Style = "Normal"

# Row 37
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E37').Style = "Normal"

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.35'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +1.30%  '
$ws.Range('E38').Style = "Normal"

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.79'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.18%  '
$ws.Range('E39').Style = "Normal"

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.85'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +4.37%  '
$ws.Range('E40').Style = "Normal"

# Row 41
$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.53'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +24.61%  '
$ws.Range('E41').Style = "Normal"

# Row 42
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.95'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('E42').Style = "Normal"

# Row 43
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0956'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +3.35%  '
$ws.Range('E43').Style = "Normal"

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '95.87'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +7.49%  '
$ws.Range('E44').Style = "Normal"

# Row 45
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.18'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +7.26%  '
$ws.Range('E45').Style = "Normal"

# Row 46
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.455.15'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('E46').Style = "Normal"

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0212'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +4.65%  '
$ws.Range('E47').Style = "Normal"

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '15.85'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +4.68%  '
$ws.Range('E48').Style = "Normal"

# Row 49
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +4.17%  '
$ws.Range('E49').Style = "Normal"

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.27'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +6.50%  '
$ws.Range('E50').Style = "Normal"

# Row 51
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.74%  '
$ws.Range('E51').Style = "Normal"
